$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure price/volume columns are treated as text so values such as
# "325.13" or "1.003" are not auto-converted to numbers by Excel,
# matching the original inline-string cell contents.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range('D2').Value = '28.843.30'
$ws.Range('E2').Value = '  +0.65%  '
$ws.Range('D3').Value = '1.876.81'
$ws.Range('E3').Value = '  -1.06%  '
$ws.Range('E4').Value = '  -0.51%  '
$ws.Range('D5').Value = '325.13'
$ws.Range('E5').Value = '  -0.56%  '
$ws.Range('D6').Value = '1.003'
$ws.Range('E6').Value = '  -0.37%  '
$ws.Range('D7').Value = '0.4601'
$ws.Range('E7').Value = '  +0.30%  '
$ws.Range('D8').Value = '0.3878'
$ws.Range('E8').Value = '  +0.18%  '
$ws.Range('D9').Value = '0.07866'
$ws.Range('E9').Value = '  -0.24%  '
$ws.Range('D10').Value = '0.9843'
$ws.Range('E10').Value = '  -1.80%  '
$ws.Range('D11').Value = '21.75'
$ws.Range('E11').Value = '  +0.04%  '
$ws.Range('D12').Value = '1.863.89'
$ws.Range('E12').Value = '  -2.30%  '
$ws.Range('D13').Value = '7.004'
$ws.Range('E13').Value = '  -1.03%  '
$ws.Range('D14').Value = '5.661'
$ws.Range('E14').Value = '  -1.01%  '
$ws.Range('D15').Value = '0.06949'
$ws.Range('E15').Value = '  -0.20%  '
$ws.Range('D16').Value = '88.25'
$ws.Range('E16').Value = '  +0.93%  '
$ws.Range('D18').Value = '0.000009973'
$ws.Range('D19').Value = '16.93'
$ws.Range('E19').Value = '  -1.28%  '
$ws.Range('E20').Value = '  -0.44%  '
$ws.Range('D21').Value = '28.851.13'
$ws.Range('E21').Value = '  +0.58%  '
$ws.Range('D22').Value = '5.262'
$ws.Range('E22').Value = '  -1.03%  '
$ws.Range('D23').Value = '10.93'
$ws.Range('E23').Value = '  -0.86%  '
$ws.Range('D24').Value = '2.084'
$ws.Range('E24').Value = '  +1.13%  '
$ws.Range('D25').Value = '155.67'
$ws.Range('E25').Value = '  +0.52%  '
$ws.Range('D26').Value = '19.30'
$ws.Range('E26').Value = '  -0.25%  '
$ws.Range('D27').Value = '5.979'
$ws.Range('E27').Value = '  +2.22%  '
$ws.Range('E28').Value = '  +0.59%  '
$ws.Range('D29').Value = '117.45'
$ws.Range('E29').Value = '  -0.99%  '
$ws.Range('D30').Value = '0.09345'
$ws.Range('E30').Value = '  +0.28%  '
$ws.Range('D31').Value = '0.9045'
$ws.Range('E31').Value = '  -2.38%  '
$ws.Range('D32').Value = '5.265'
$ws.Range('E32').Value = '  -0.91%  '
$ws.Range('E33').Value = '  -0.84%  '
$ws.Range('D34').Value = '3.266'
$ws.Range('E34').Value = '  -0.07%  '
$ws.Range('D35').Value = '1.188'
$ws.Range('E35').Value = '  +2.43%  '
$ws.Range('D36').Value = '0.05767'
$ws.Range('E36').Value = '  +0.30%  '
$ws.Range('E37').Value = '  +0.04%  '
$ws.Range('E38').Value = '  -0.40%  '
$ws.Range('D39').Value = '7.681'
$ws.Range('E39').Value = '  -0.86%  '
$ws.Range('D40').Value = '0.5657'
$ws.Range('E40').Value = '  +0.41%  '
$ws.Range('D41').Value = '0.1767'
$ws.Range('E41').Value = '  -0.96%  '
$ws.Range('D42').Value = '9.677'
$ws.Range('E42').Value = '  -1.04%  '
$ws.Range('D43').Value = '2.255'
$ws.Range('E43').Value = '  +2.28%  '
$ws.Range('D44').Value = '11.89'
$ws.Range('E44').Value = '  +0.89%  '
$ws.Range('D45').Value = '0.5350'
$ws.Range('E45').Value = '  +0.29%  '
$ws.Range('D46').Value = '0.07036'
$ws.Range('E46').Value = '  -1.87%  '
$ws.Range('D47').Value = '1.845'
$ws.Range('E47').Value = '  +0.58%  '
$ws.Range('D48').Value = '113.06'
$ws.Range('E48').Value = '  +0.35%  '
$ws.Range('D49').Value = '2.506'
$ws.Range('E49').Value = '  +1.89%  '
$ws.Range('E50').Value = '  -5.18%  '
$ws.Range('D51').Value = '70.69'
$ws.Range('E51').Value = '  +0.03%  '

# Restore the default (Normal) cell style so no stray number-format
# styling is left behind on the edited range.
$ws.Range("D2:E51").Style = "Normal"

